$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must remain literal text
# (matching the original inline-string cell type) rather than being auto-
# converted to numbers, so force text format for the whole column up front,
# assign the new values, then restore normal styling.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '68.765.20'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.861.62'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '600.30'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').Value = '162.16'
$ws.Range('E6').Value = '  -2.71%  '
$ws.Range('D7').Value = '3.859.60'
$ws.Range('E7').Value = '  +3.00%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('D10').Value = '0.167'
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').Value = '36.84'
$ws.Range('E13').Value = '  -2.80%  '
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').Value = '4.509.43'
$ws.Range('E15').Value = '  +3.12%  '
$ws.Range('D16').Value = '3.836.98'
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('D17').Value = '68.949.61'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '7.56'
$ws.Range('E18').Value = '  +2.84%  '
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '17.12'
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '11.35'
$ws.Range('E21').Value = '  +2.91%  '
$ws.Range('D22').Value = '483.45'
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('E24').Value = '  +6.47%  '
$ws.Range('D25').Value = '83.90'
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('D27').Value = '12.07'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').Value = '4.015.22'
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('D32').Value = '7.85'
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').Value = '32.19'
$ws.Range('E33').Value = '  +2.47%  '
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('D35').Value = '3.809.67'
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('D37').Value = '1.04'
$ws.Range('E37').Value = '  +1.99%  '
$ws.Range('E38').Value = '  +1.42%  '
$ws.Range('D39').Value = '5.87'
$ws.Range('E39').Value = '  -1.35%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = '0.318'
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('D42').Value = '2.96'
$ws.Range('E42').Value = '  -2.18%  '
$ws.Range('D43').Value = '431.31'
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('D44').Value = '48.50'
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('D45').Value = '1.97'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('E47').Value = '  -0.99%  '
$ws.Range('D48').Value = '143.78'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('D49').Value = '2.840.99'
$ws.Range('E49').Value = '  +1.91%  '
$ws.Range('D50').Value = '0.0358'
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('D51').Value = '25.82'
$ws.Range('E51').Value = '  +12.98%  '

$ws.Range("D2:D51").Style = "Normal"
